# Adds a task which opens a url instead of a text screen or a command.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "url" column header (column F)
$ws.Range("F1").Value = "url"

# New task row (row 4): task_name = "url", duration in minutes = 3
$ws.Range("A4").Value = "url"
$ws.Range("B4").Value = 3

# The url task stores its target address as a real hyperlink in column F,
# displaying the url text and using the built-in "Hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.google.com", "", "", "https://www.google.com")

# Nudge & grow the "HowTo" textbox so the new url instructions still fit.
$shape = $ws.Shapes.Item(1)
$shape.Left = $shape.Left + 36
$shape.Top = $shape.Top + 1
$shape.Height = $shape.Height + 40
